$wb = $excel.ActiveWorkbook

# --- "ingredients" sheet: move selection from B6 to B4 ---
$wsIngredients = $wb.Worksheets.Item("ingredients")
$wsIngredients.Range("B4").Select()

# --- "ingredient-recipemap" sheet: re-select A11 so the scrolled
#     topLeftCell="A22" view state resets back to showing the selection ---
$wsMap = $wb.Worksheets.Item("ingredient-recipemap")
$wsMap.Range("A11").Select()

# --- "ingredients table" sheet: add a datatype column (C) of "int" for
#     every ingredient row (2-14), then update the selection ---
$wsTable = $wb.Worksheets.Item("ingredients table")
for ($r = 2; $r -le 14; $r++) {
    $wsTable.Cells.Item($r, 3).Value = "int"
}
$wsTable.Range("B2:C14").Select()

Write-Output "done"
